$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update time_taken (F column) timestamps on the "data" sheet ---
$ws.Range("F2").Value = "2021-10-05 14:35:34.242499"
$ws.Range("F3").Value = "2021-10-05 14:35:34.242507"
$ws.Range("F4").Value = "2021-10-05 14:35:34.242510"
$ws.Range("F5").Value = "2021-10-05 14:35:34.242512"
$ws.Range("F6").Value = "2021-10-05 14:35:34.242515"
$ws.Range("F7").Value = "2021-10-05 14:35:34.242518"
$ws.Range("F8").Value = "2021-10-05 14:35:34.242521"
$ws.Range("F9").Value = "2021-10-05 14:35:34.242523"
$ws.Range("F10").Value = "2021-10-05 14:35:34.242526"
$ws.Range("F11").Value = "2021-10-05 14:35:34.242529"
$ws.Range("F12").Value = "2021-10-05 14:35:34.242531"
$ws.Range("F13").Value = "2021-10-05 14:35:34.242534"
$ws.Range("F14").Value = "2021-10-05 14:35:34.242536"
$ws.Range("F15").Value = "2021-10-05 14:35:34.242538"
$ws.Range("F16").Value = "2021-10-05 14:35:34.242541"
$ws.Range("F17").Value = "2021-10-05 14:35:34.242543"
$ws.Range("F18").Value = "2021-10-05 14:35:34.242546"
$ws.Range("F19").Value = "2021-10-05 14:35:34.242549"
$ws.Range("F20").Value = "2021-10-05 14:35:34.242551"
$ws.Range("F21").Value = "2021-10-05 14:35:34.242553"
$ws.Range("F22").Value = "2021-10-05 14:35:34.242556"
$ws.Range("F23").Value = "2021-10-05 14:35:34.242558"
$ws.Range("F24").Value = "2021-10-05 14:35:34.242561"
$ws.Range("F25").Value = "2021-10-05 14:35:34.242563"
$ws.Range("F26").Value = "2021-10-05 14:35:34.242566"
$ws.Range("F27").Value = "2021-10-05 14:35:34.242568"
$ws.Range("F28").Value = "2021-10-05 14:35:34.242571"
$ws.Range("F29").Value = "2021-10-05 14:35:34.242573"
$ws.Range("F30").Value = "2021-10-05 14:35:34.242576"
$ws.Range("F31").Value = "2021-10-05 14:35:34.242578"
$ws.Range("F32").Value = "2021-10-05 14:35:34.242580"
$ws.Range("F33").Value = "2021-10-05 14:35:34.242583"
$ws.Range("F34").Value = "2021-10-05 14:35:34.242586"
$ws.Range("F35").Value = "2021-10-05 14:35:34.242588"
$ws.Range("F36").Value = "2021-10-05 14:35:34.242591"
$ws.Range("F37").Value = "2021-10-05 14:35:34.242593"
$ws.Range("F38").Value = "2021-10-05 14:35:34.242596"
$ws.Range("F39").Value = "2021-10-05 14:35:34.242598"
$ws.Range("F40").Value = "2021-10-05 14:35:34.242600"
$ws.Range("F41").Value = "2021-10-05 14:35:34.242603"
$ws.Range("F42").Value = "2021-10-05 14:35:34.242606"
$ws.Range("F43").Value = "2021-10-05 14:35:34.242608"
$ws.Range("F44").Value = "2021-10-05 14:35:34.242611"
$ws.Range("F45").Value = "2021-10-05 14:35:34.242613"
$ws.Range("F46").Value = "2021-10-05 14:35:34.242616"
$ws.Range("F47").Value = "2021-10-05 14:35:34.242618"
$ws.Range("F48").Value = "2021-10-05 14:35:34.242621"
$ws.Range("F49").Value = "2021-10-05 14:35:34.242623"
$ws.Range("F50").Value = "2021-10-05 14:35:34.242625"
$ws.Range("F51").Value = "2021-10-05 14:35:34.242628"
$ws.Range("F52").Value = "2021-10-05 14:35:34.242631"
$ws.Range("F53").Value = "2021-10-05 14:35:34.242633"
$ws.Range("F54").Value = "2021-10-05 14:35:34.242636"
$ws.Range("F55").Value = "2021-10-05 14:35:34.242639"
$ws.Range("F56").Value = "2021-10-05 14:35:34.242641"
$ws.Range("F57").Value = "2021-10-05 14:35:34.242643"
$ws.Range("F58").Value = "2021-10-05 14:35:34.242646"
$ws.Range("F59").Value = "2021-10-05 14:35:34.242648"
$ws.Range("F60").Value = "2021-10-05 14:35:34.242651"
$ws.Range("F61").Value = "2021-10-05 14:35:34.242653"
$ws.Range("F62").Value = "2021-10-05 14:35:34.242656"
$ws.Range("F63").Value = "2021-10-05 14:35:34.242658"
$ws.Range("F64").Value = "2021-10-05 14:35:34.242661"
$ws.Range("F65").Value = "2021-10-05 14:35:34.242663"
$ws.Range("F66").Value = "2021-10-05 14:35:34.242667"
$ws.Range("F67").Value = "2021-10-05 14:35:34.242669"
$ws.Range("F68").Value = "2021-10-05 14:35:34.242672"
$ws.Range("F69").Value = "2021-10-05 14:35:34.242674"
$ws.Range("F70").Value = "2021-10-05 14:35:34.242677"
$ws.Range("F71").Value = "2021-10-05 14:35:34.242679"
$ws.Range("F72").Value = "2021-10-05 14:35:34.242682"
$ws.Range("F73").Value = "2021-10-05 14:35:34.242684"
$ws.Range("F74").Value = "2021-10-05 14:35:34.242686"
$ws.Range("F75").Value = "2021-10-05 14:35:34.242689"
$ws.Range("F76").Value = "2021-10-05 14:35:34.242691"
$ws.Range("F77").Value = "2021-10-05 14:35:34.242694"
$ws.Range("F78").Value = "2021-10-05 14:35:34.242698"
$ws.Range("F79").Value = "2021-10-05 14:35:34.242701"
$ws.Range("F80").Value = "2021-10-05 14:35:34.242704"
$ws.Range("F81").Value = "2021-10-05 14:35:34.242706"
$ws.Range("F82").Value = "2021-10-05 14:35:34.242709"
$ws.Range("F83").Value = "2021-10-05 14:35:34.242711"
$ws.Range("F84").Value = "2021-10-05 14:35:34.242714"
$ws.Range("F85").Value = "2021-10-05 14:35:34.242716"
$ws.Range("F86").Value = "2021-10-05 14:35:34.242719"
$ws.Range("F87").Value = "2021-10-05 14:35:34.242721"
$ws.Range("F88").Value = "2021-10-05 14:35:34.242724"
$ws.Range("F89").Value = "2021-10-05 14:35:34.242727"
$ws.Range("F90").Value = "2021-10-05 14:35:34.242729"
$ws.Range("F91").Value = "2021-10-05 14:35:34.242732"
$ws.Range("F92").Value = "2021-10-05 14:35:34.242735"
$ws.Range("F93").Value = "2021-10-05 14:35:34.242737"
$ws.Range("F94").Value = "2021-10-05 14:35:34.242741"
$ws.Range("F95").Value = "2021-10-05 14:35:34.242744"
$ws.Range("F96").Value = "2021-10-05 14:35:34.242746"
$ws.Range("F97").Value = "2021-10-05 14:35:34.242749"
$ws.Range("F98").Value = "2021-10-05 14:35:34.242752"

# --- Add the new "metadata" sheet after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newWs.Name = "metadata"

# Header row (row 1)
$newWs.Range("B1").Value = "data_name"
$newWs.Range("C1").Value = "data_id"
$newWs.Range("D1").Value = "data_version"
$newWs.Range("E1").Value = "data_version_created"
$newWs.Range("F1").Value = "panel_query_time"
$newWs.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "Renal Glomerular Disease_SuperPanel"
$newWs.Range("C2").Value = 262

# "1.20" must stay literal text (not be normalised to the number 1.2) while the
# cell itself keeps the sheet's plain default style (no explicit style index).
# Stage the text in a scratch cell that is formatted as Text, then copy only
# the VALUE (not the format) into D2, so D2 ends up styleless but still holds
# the string "1.20". The scratch cell is cleared (and its format wiped) after.
$newWs.Range("Z1").NumberFormat = "@"
$newWs.Range("Z1").Value = "1.20"
$newWs.Range("Z1").Copy()
$newWs.Range("D2").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$newWs.Range("Z1").Clear()

$newWs.Range("E2").Value = "2021-08-06T01:08:30.713247Z"
$newWs.Range("F2").Value = "2021-10-05 14:35:34.239076"
$newWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/262/?format=json"

# Re-use the bold/bordered/centered header style already used on the "data" sheet
# (style index 1) by copy/paste-special of formats only - avoids creating new style entries.
$ws.Range("B1:F1").Copy()
$newWs.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$newWs.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave selection at A1 on the new sheet, then re-activate the "data" sheet to match original state
$newWs.Range("A1").Select()
$ws.Activate()
